# Revert "Merge branch 'cluster-algorithm-sunflower' ..." -- restores the
# pre-merge scenario-parameter/specs-data values and UI view state.

$wb = $excel.ActiveWorkbook

# --- ScenarioParameters: tweak two scenario parameter values -------------
$wsParams = $wb.Worksheets.Item("ScenarioParameters")
$wsParams.Range("H2").Value = 1.2
$wsParams.Range("F3").Value = 0.084

# --- SpecsData: revert the two investment-cost figures --------------------
$wsSpecs = $wb.Worksheets.Item("SpecsData")
$wsSpecs.Range("I2").Value = 30000000
$wsSpecs.Range("J2").Value = 26600000

# --- SpecsDataCalib: same two figures --------------------------------------
$wsCalib = $wb.Worksheets.Item("SpecsDataCalib")
$wsCalib.Range("I2").Value = 30000000
$wsCalib.Range("J2").Value = 26600000

# --- SpecsDataCalib1: same two figures -------------------------------------
$wsCalib1 = $wb.Worksheets.Item("SpecsDataCalib1")
$wsCalib1.Range("I2").Value = 30000000
$wsCalib1.Range("J2").Value = 26600000

# --- Restore each sheet's own last-used selection --------------------------
$wsParams.Activate()
$wsParams.Range("I17").Select()

$wsCalib.Activate()
$wsCalib.Range("Y16").Select()

$wsCalib1.Activate()
$wsCalib1.Range("O8").Select()

# SpecsData becomes the active (tabSelected) sheet again, as the last one
# activated, with its own restored selection/scroll position.
$wsSpecs.Activate()
$wsSpecs.Range("Q2").Select()
